$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold / bordered / centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2-29: column I (I0) and column J (IF)
$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 6; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 4; 21 = 1; 22 = 1; 23 = 4; 24 = 5; 25 = 1; 26 = 1; 27 = 3; 28 = 2; 29 = 1
}
$jValues = @{
    2 = 4; 3 = 5; 4 = 7; 5 = 6; 6 = 6; 7 = 7; 8 = 7; 9 = 5; 10 = 5;
    11 = 4; 12 = 5; 13 = 8; 14 = 5; 15 = 6; 16 = 7; 17 = 6; 18 = 6; 19 = 5;
    20 = 7; 21 = 3; 22 = 6; 23 = 7; 24 = 6; 25 = 6; 26 = 5; 27 = 4; 28 = 3; 29 = 1
}

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
